$d = $word.ActiveDocument

# wdReplaceAll = 2
$replaceAll = 2

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, $replaceAll)
}

Replace-All "0.9813 and coefficients" "0.9805 and coefficients"
Replace-All "2.9841, and" "3.018, and"
Replace-All "2.0021" "1.999"

Replace-All "0.9807 and coefficients" "0.98 and coefficients"
Replace-All "2.9832, and" "3.0169, and"
Replace-All "2.0014" "1.9986"

Replace-All "array([0.99184944]) and coefficients" "array([0.98111268]) and coefficients"
Replace-All "array([2.95951084]), and" "array([3.00982231]), and"
Replace-All "array([2.00995977])" "array([1.96969333])"
